$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = 0.2622798531249657
$ws.Cells.Item(5, 2).Value = 0.4989385597712612
$ws.Cells.Item(6, 2).Value = 0.1216139147831967
$ws.Cells.Item(7, 2).Value = 0.6548011357729422
$ws.Cells.Item(8, 2).Value = 0.9925931342935428
$ws.Cells.Item(9, 2).Value = 0.4197174453273789
$ws.Cells.Item(10, 2).Value = 0.4423121297183457
$ws.Cells.Item(11, 2).Value = 0.3248472793844002
$ws.Cells.Item(12, 2).Value = 0.07907397708568453
$ws.Cells.Item(13, 2).Value = 0.2409202938150372
$ws.Cells.Item(14, 2).Value = 0.1901876995395071
$ws.Cells.Item(15, 2).Value = 0.3234244561263893
$ws.Cells.Item(16, 2).Value = 0.5460928092329453
$ws.Cells.Item(17, 2).Value = 0.24993599858945
$ws.Cells.Item(18, 2).Value = -0.0008109998315291023
$ws.Cells.Item(19, 2).Value = 0.2240000000000038
$ws.Cells.Item(20, 2).Value = 0.2694296772580742
$ws.Cells.Item(21, 2).Value = 0.3037528065631818
$ws.Cells.Item(22, 2).Value = 0.102683801342323
$ws.Cells.Item(23, 2).Value = 0.4016405696262808
$ws.Cells.Item(24, 2).Value = 0.4670353440379102
$ws.Cells.Item(25, 2).Value = 0.3229903499871796
$ws.Cells.Item(26, 2).Value = 0.344825264879205
$ws.Cells.Item(27, 2).Value = 0.2849999999999966
$ws.Cells.Item(28, 2).Value = 0.2296800116567681
$ws.Cells.Item(29, 2).Value = 0.319136520886417
$ws.Cells.Item(30, 2).Value = 0.3192030404897963
$ws.Cells.Item(31, 2).Value = 0.3838834518010543
$ws.Cells.Item(32, 2).Value = 0.1721310865893741
$ws.Cells.Item(33, 2).Value = 0.3124151146021887
$ws.Cells.Item(34, 2).Value = 1.338781480463297
$ws.Cells.Item(35, 2).Value = 0.6782258938801533
$ws.Cells.Item(36, 2).Value = 0.5169063711852431
$ws.Cells.Item(37, 2).Value = 0.5229245835751897
$ws.Cells.Item(38, 2).Value = 0.6999999999999886
$ws.Cells.Item(39, 2).Value = 0.8499999999999943
$ws.Cells.Item(40, 2).Value = 0.4000000000000057
$ws.Cells.Item(41, 2).Value = 0.3500000000000085
$ws.Cells.Item(42, 2).Value = 0.1700920104450461
$ws.Cells.Item(43, 2).Value = 0.2000000000000028
$ws.Cells.Item(44, 2).Value = 0.4999999999999858
$ws.Cells.Item(45, 2).Value = 0.5999999999999943
$ws.Cells.Item(46, 2).Value = 0.5
$ws.Cells.Item(47, 2).Value = 0.6648928967370864
$ws.Cells.Item(48, 2).Value = 0.7
$ws.Cells.Item(49, 2).Value = 0.3999908011606834
$ws.Cells.Item(50, 2).Value = 0.5200000000000102
$ws.Cells.Item(51, 2).Value = 0.5999999999999943
$ws.Cells.Item(52, 2).Value = 0.5
$ws.Cells.Item(53, 2).Value = 0.3499999990000049
$ws.Cells.Item(54, 2).Value = 0.7000000000000171
$ws.Cells.Item(55, 2).Value = 0.7000000000000171
$ws.Cells.Item(56, 2).Value = 0.7999999999999972
$ws.Cells.Item(57, 2).Value = -0.7999999999999972
$ws.Cells.Item(58, 2).Value = 0.4999999999999858
$ws.Cells.Item(59, 2).Value = [double]"2.842170943040401E-14"
$ws.Cells.Item(60, 2).Value = 0.5000000000000142
$ws.Cells.Item(61, 2).Value = 0.7000000000000171
$ws.Cells.Item(62, 2).Value = -0.9999999999999858
$ws.Cells.Item(63, 2).Value = -0.5
$ws.Cells.Item(64, 2).Value = 2.799999999999997
$ws.Cells.Item(65, 2).Value = 0
$ws.Cells.Item(66, 2).Value = 0.2000000000000028
$ws.Cells.Item(67, 2).Value = -0.09999999999999432
$ws.Cells.Item(68, 2).Value = 0.2000000000000028
$ws.Cells.Item(69, 2).Value = 0.09999999999999432
$ws.Cells.Item(70, 2).Value = 0.09999999999999432
$ws.Cells.Item(71, 2).Value = 0
$ws.Cells.Item(72, 2).Value = 0.4999999999999858
$ws.Cells.Item(73, 2).Value = 0.4000000000000199

# Delete rows 74 to 82 (9 rows) - shift cells up
$ws.Rows("74:82").Delete()
